# Add the new weekly ranking sheet for 2025-12-17, placed after 2025-12-10,
# mirroring the layout/styling of the existing weekly sheets.
$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2025-12-17"

# --- Header row (rank / title / volume / publisher) ---
$ws.Cells.Item(1, 1).Value = 'rank'
$ws.Cells.Item(1, 2).Value = 'title'
$ws.Cells.Item(1, 3).Value = 'volume'
$ws.Cells.Item(1, 4).Value = 'publisher'

$headerRange = $ws.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1         # xlContinuous
$headerRange.Borders.Weight = 2            # xlThin

# --- Data rows ---
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 'ONE PIECE'
$ws.Cells.Item(2, 3).Value = 113
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 'ワールドトリガー'
$ws.Cells.Item(3, 3).Value = 29
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 'SAKAMOTO DAYS'
$ws.Cells.Item(4, 3).Value = 25
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 'ザ・ファブル The third secret'
$ws.Cells.Item(5, 3).Value = 3
$ws.Cells.Item(5, 3).Interior.Color = 13499135
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = '青の祓魔師'
$ws.Cells.Item(6, 3).Value = 33
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 'アオのハコ'
$ws.Cells.Item(7, 3).Value = 23
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = '薬屋のひとりごと'
$ws.Cells.Item(8, 3).Value = 16
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = '魔入りました!入間くん'
$ws.Cells.Item(9, 3).Value = 46
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = '逃げ上手の若君'
$ws.Cells.Item(10, 3).Value = 23
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = '信じていた仲間達にダンジョン奥地で殺されかけたがギフト『無限ガチャ』でレベル9999の仲間達を手に入れて元パーティーメンバーと世界に復讐&『ざまぁ!』します!'
$ws.Cells.Item(11, 3).Value = 20
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = 'アマチュアビジランテ'
$ws.Cells.Item(12, 3).Value = 6
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = '薫る花は凛と咲く'
$ws.Cells.Item(13, 3).Value = 21
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = '魔男のイチ'
$ws.Cells.Item(14, 3).Value = 6
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = '勇者パーティを追い出された器用貧乏 ~パーティ事情で付与術士をやっていた剣士、万能へと至る~'
$ws.Cells.Item(15, 3).Value = 17
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = 'ブルーロック'
$ws.Cells.Item(16, 3).Value = 36
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = '怪物事変'
$ws.Cells.Item(17, 3).Value = 24
$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = 'ウィッチウォッチ'
$ws.Cells.Item(18, 3).Value = 24
$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = 'ブルーピリオド'
$ws.Cells.Item(19, 3).Value = 18
$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = 'グラぱらっ!'
$ws.Cells.Item(20, 3).Value = 10
$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = '桃源暗鬼'
$ws.Cells.Item(21, 3).Value = 27
$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = 'ねずみの初恋'
$ws.Cells.Item(22, 3).Value = 8
$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = 'ワンパンマン'
$ws.Cells.Item(23, 3).Value = 35
$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 2).Value = 'ラーメン赤猫'
$ws.Cells.Item(24, 3).Value = 13
$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 2).Value = '味方が弱すぎて補助魔法に徹していた宮廷魔法師、追放されて最強を目指す'
$ws.Cells.Item(25, 3).Value = 18
$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(26, 2).Value = '似鳥教授の可愛い悪あがき'
$ws.Cells.Item(26, 3).Value = 2
$ws.Cells.Item(26, 3).Interior.Color = 13499135
$ws.Cells.Item(27, 1).Value = 26
$ws.Cells.Item(27, 2).Value = 'アザミヤコを好きになる'
$ws.Cells.Item(27, 3).Value = 1
$ws.Cells.Item(27, 3).Interior.Color = 13499135
$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = '神血の救世主~0.00000001%を引き当て最強へ~'
$ws.Cells.Item(28, 3).Value = 7
$ws.Cells.Item(29, 1).Value = 28
$ws.Cells.Item(29, 2).Value = '実は俺、最強でした?'
$ws.Cells.Item(29, 3).Value = 18
$ws.Cells.Item(30, 1).Value = 29
$ws.Cells.Item(30, 2).Value = 'WIND BREAKER'
$ws.Cells.Item(30, 3).Value = 24
$ws.Cells.Item(31, 1).Value = 30
$ws.Cells.Item(31, 2).Value = 'FAIRY TAIL 100 YEARS QUEST'
$ws.Cells.Item(31, 3).Value = 22
$ws.Cells.Item(32, 1).Value = 31
$ws.Cells.Item(32, 2).Value = 'いきなり結婚宣言~裏切られ絶望した私に待っていたのは溺愛でした~1'
$ws.Cells.Item(32, 3).Value = 1
$ws.Cells.Item(32, 3).Interior.Color = 13499135
$ws.Cells.Item(33, 1).Value = 32
$ws.Cells.Item(33, 2).Value = '鬼上司のヤキモチが可愛すぎます!!'
$ws.Cells.Item(33, 3).Value = 3
$ws.Cells.Item(33, 3).Interior.Color = 13499135
$ws.Cells.Item(34, 1).Value = 33
$ws.Cells.Item(34, 2).Value = '私と上司の内緒の事情'
$ws.Cells.Item(34, 3).Value = 15
$ws.Cells.Item(35, 1).Value = 34
$ws.Cells.Item(35, 2).Value = 'その悪役貴族、ママヒロインが好きすぎる'
$ws.Cells.Item(35, 3).Value = 1
$ws.Cells.Item(35, 3).Interior.Color = 13499135
$ws.Cells.Item(36, 1).Value = 35
$ws.Cells.Item(36, 2).Value = '嘘つきカノジョの影盛さん(フルカラー)'
$ws.Cells.Item(36, 3).Value = 1
$ws.Cells.Item(36, 3).Interior.Color = 13499135
$ws.Cells.Item(37, 1).Value = 36
$ws.Cells.Item(37, 2).Value = 'ダンジョンに置き去りにされたので定住した~そしてダンジョンマスターへ~'
$ws.Cells.Item(37, 3).Value = 5
$ws.Cells.Item(38, 1).Value = 37
$ws.Cells.Item(38, 2).Value = '転生したらスライムだった件'
$ws.Cells.Item(38, 3).Value = 30
$ws.Cells.Item(39, 1).Value = 38
$ws.Cells.Item(39, 2).Value = 'いつでも自宅に帰れる俺は、異世界で行商人をはじめました'
$ws.Cells.Item(39, 3).Value = 7
$ws.Cells.Item(40, 1).Value = 39
$ws.Cells.Item(40, 2).Value = '異世界領地改革~土魔法で始める公共事業~'
$ws.Cells.Item(40, 3).Value = 8
$ws.Cells.Item(41, 1).Value = 40
$ws.Cells.Item(41, 2).Value = 'ふつうの軽音部'
$ws.Cells.Item(41, 3).Value = 9
$ws.Cells.Item(42, 1).Value = 41
$ws.Cells.Item(42, 2).Value = 'SPY×FAMILY'
$ws.Cells.Item(42, 3).Value = 16
$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = 'なかなか稀少な光谷さん'
$ws.Cells.Item(43, 3).Value = 1
$ws.Cells.Item(43, 3).Interior.Color = 13499135
$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = '魔法歌姫マジカルギンガ 第24話'
$ws.Cells.Item(44, 3).Value = 24
$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = 'ハニトラ・ホームステイ~北欧人妻アンナさんの危険な香り~'
$ws.Cells.Item(45, 3).Value = 1
$ws.Cells.Item(45, 3).Interior.Color = 13499135
$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = '神血の救世主~0.00000001%を引き当て最強へ~'
$ws.Cells.Item(46, 3).Value = 8
$ws.Cells.Item(47, 1).Value = 46
$ws.Cells.Item(47, 2).Value = '転生したら第七王子だったので、気ままに魔術を極めます'
$ws.Cells.Item(47, 3).Value = 21
$ws.Cells.Item(48, 1).Value = 47
$ws.Cells.Item(48, 2).Value = 'この冒険者、人類史最強です~外れスキル『鑑定』が『継承』に覚醒したので、数多の英雄たちの力を受け継ぎ無双する~'
$ws.Cells.Item(48, 3).Value = 6
$ws.Cells.Item(49, 1).Value = 48
$ws.Cells.Item(49, 2).Value = 'こじらせ令嬢ですが、実は王子だった男友達に溺愛されています第1話'
$ws.Cells.Item(49, 3).Value = 1
$ws.Cells.Item(49, 3).Interior.Color = 13499135
$ws.Cells.Item(50, 1).Value = 49
$ws.Cells.Item(50, 2).Value = 'うたの☆プリンセスさまっ♪BACK to the IDOL1話'
$ws.Cells.Item(50, 3).Value = 1
$ws.Cells.Item(50, 3).Interior.Color = 13499135
$ws.Cells.Item(51, 1).Value = 50
$ws.Cells.Item(51, 2).Value = 'モラハラ不倫~私はやっぱり欠陥品1'
$ws.Cells.Item(51, 3).Value = 1
$ws.Cells.Item(51, 3).Interior.Color = 13499135
$ws.Cells.Item(52, 1).Value = 51
$ws.Cells.Item(52, 2).Value = '悪役令嬢は全力でグータラしたいのに、隣国皇太子が溺愛してくる。なぜ。 (エンジェライトコミックス)'
$ws.Cells.Item(52, 3).Value = 1
$ws.Cells.Item(52, 3).Interior.Color = 13499135
$ws.Cells.Item(53, 1).Value = 52
$ws.Cells.Item(53, 2).Value = '別れ際、アイツは私にキスをする'
$ws.Cells.Item(53, 3).Value = 3
$ws.Cells.Item(53, 3).Interior.Color = 13499135
$ws.Cells.Item(54, 1).Value = 53
$ws.Cells.Item(54, 2).Value = '鬼上司のヤキモチが可愛すぎます!!'
$ws.Cells.Item(54, 3).Value = 8
$ws.Cells.Item(55, 1).Value = 54
$ws.Cells.Item(55, 2).Value = '私と上司の内緒の事情'
$ws.Cells.Item(55, 3).Value = 12
$ws.Cells.Item(56, 1).Value = 55
$ws.Cells.Item(56, 2).Value = '私と上司の内緒の事情'
$ws.Cells.Item(56, 3).Value = 13
$ws.Cells.Item(57, 1).Value = 56
$ws.Cells.Item(57, 2).Value = '私と上司の内緒の事情'
$ws.Cells.Item(57, 3).Value = 14
$ws.Cells.Item(58, 1).Value = 57
$ws.Cells.Item(58, 2).Value = 'なかなか稀少な光谷さん'
$ws.Cells.Item(58, 3).Value = 2
$ws.Cells.Item(58, 3).Interior.Color = 13499135
$ws.Cells.Item(59, 1).Value = 58
$ws.Cells.Item(59, 2).Value = 'なかなか稀少な光谷さん'
$ws.Cells.Item(59, 3).Value = 3
$ws.Cells.Item(59, 3).Interior.Color = 13499135
$ws.Cells.Item(60, 1).Value = 59
$ws.Cells.Item(60, 2).Value = 'アザミヤコを好きになる'
$ws.Cells.Item(60, 3).Value = 2
$ws.Cells.Item(60, 3).Interior.Color = 13499135
$ws.Cells.Item(61, 1).Value = 60
$ws.Cells.Item(61, 2).Value = 'アザミヤコを好きになる'
$ws.Cells.Item(61, 3).Value = 3
$ws.Cells.Item(61, 3).Interior.Color = 13499135
$ws.Cells.Item(62, 1).Value = 61
$ws.Cells.Item(62, 2).Value = 'その悪役貴族、ママヒロインが好きすぎる'
$ws.Cells.Item(62, 3).Value = 2
$ws.Cells.Item(62, 3).Interior.Color = 13499135
$ws.Cells.Item(63, 1).Value = 62
$ws.Cells.Item(63, 2).Value = 'その悪役貴族、ママヒロインが好きすぎる'
$ws.Cells.Item(63, 3).Value = 3
$ws.Cells.Item(63, 3).Interior.Color = 13499135
$ws.Cells.Item(64, 1).Value = 63
$ws.Cells.Item(64, 2).Value = '灰かぶりの天使'
$ws.Cells.Item(64, 3).Value = 1
$ws.Cells.Item(64, 3).Interior.Color = 13499135
$ws.Cells.Item(65, 1).Value = 64
$ws.Cells.Item(65, 2).Value = '死霊術師ウェルツの平和論 WEBコミックガンマぷらす連載版 第1話'
$ws.Cells.Item(65, 3).Value = 1
$ws.Cells.Item(65, 3).Interior.Color = 13499135
$ws.Cells.Item(66, 1).Value = 65
$ws.Cells.Item(66, 2).Value = '聖なる加護持ち令嬢は、騎士を目指しているので聖女にはなりません。 WEBコミックガンマぷらす連載版 第一話'
$ws.Cells.Item(66, 3).Value = 1
$ws.Cells.Item(66, 3).Interior.Color = 13499135
$ws.Cells.Item(67, 1).Value = 66
$ws.Cells.Item(67, 2).Value = '転生したらポンコツメイドと呼ばれていました'
$ws.Cells.Item(67, 3).Value = 1
$ws.Cells.Item(67, 3).Interior.Color = 13499135
$ws.Cells.Item(68, 1).Value = 67
$ws.Cells.Item(68, 2).Value = '聖人公爵様がラスボスだということを私だけが知っている ストーリアダッシュ連載版 第1話'
$ws.Cells.Item(68, 3).Value = 1
$ws.Cells.Item(68, 3).Interior.Color = 13499135
$ws.Cells.Item(69, 1).Value = 68
$ws.Cells.Item(69, 2).Value = '義妹に婚約者を奪われたので、好きに生きようと思います。 ストーリアダッシュ連載版 第1話'
$ws.Cells.Item(69, 3).Value = 1
$ws.Cells.Item(69, 3).Interior.Color = 13499135
$ws.Cells.Item(70, 1).Value = 69
$ws.Cells.Item(70, 2).Value = '幽閉令嬢の気ままな異世界生活~転生ライフを楽しんでいるので、邪魔しに来ないでくれませんか、元婚約者様?~ ストーリアダッシュ連載版 第1話'
$ws.Cells.Item(70, 3).Value = 1
$ws.Cells.Item(70, 3).Interior.Color = 13499135
$ws.Cells.Item(71, 1).Value = 70
$ws.Cells.Item(71, 2).Value = '透明人間になったわたしと、わたしに興味がない(はずの)夫の奇妙な三か月間 ストーリアダッシュ連載版 第1話'
$ws.Cells.Item(71, 3).Value = 1
$ws.Cells.Item(71, 3).Interior.Color = 13499135
$ws.Cells.Item(72, 1).Value = 71
$ws.Cells.Item(72, 2).Value = '溺愛策士な護衛騎士は純粋培養令嬢に意地悪したい。 ストーリアダッシュ連載版 第1話'
$ws.Cells.Item(72, 3).Value = 1
$ws.Cells.Item(72, 3).Interior.Color = 13499135
$ws.Cells.Item(73, 1).Value = 72
$ws.Cells.Item(73, 2).Value = 'エルフさんの魔法料理店 妖精女王として転生したけれど、まずはのんびりお料理作りまくります! ストーリアダッシュ連載版 第1話'
$ws.Cells.Item(73, 3).Value = 1
$ws.Cells.Item(73, 3).Interior.Color = 13499135
$ws.Cells.Item(74, 1).Value = 73
$ws.Cells.Item(74, 2).Value = '信じていた仲間達にダンジョン奥地で殺されかけたがギフト『無限ガチャ』でレベル9999の仲間達を手に入れて元パーティーメンバーと世界に復讐&『ざまぁ!』します!'
$ws.Cells.Item(74, 3).Value = 19
$ws.Cells.Item(75, 1).Value = 74
$ws.Cells.Item(75, 2).Value = '再召喚された勇者は一般人として生きていく?'
$ws.Cells.Item(75, 3).Value = 11
$ws.Cells.Item(76, 1).Value = 75
$ws.Cells.Item(76, 2).Value = 'サンキューピッチ'
$ws.Cells.Item(76, 3).Value = 4
$ws.Cells.Item(77, 1).Value = 76
$ws.Cells.Item(77, 2).Value = '彼女、お借りします'
$ws.Cells.Item(77, 3).Value = 43
$ws.Cells.Item(78, 1).Value = 77
$ws.Cells.Item(78, 2).Value = 'わたしの幸せな結婚'
$ws.Cells.Item(78, 3).Value = 37
$ws.Cells.Item(79, 1).Value = 78
$ws.Cells.Item(79, 2).Value = '零細奴隷商人、一人も奴隷が売れなかったので売れ残り少女たちと辺境でスローライフをする~毎日優しく接していたら、いつの間にか勝手に魔物を狩るようになってきた。え、この子たち最強種の魔族だったの?~(ノヴァコミックス)2'
$ws.Cells.Item(79, 3).Value = 2
$ws.Cells.Item(79, 3).Interior.Color = 13499135
$ws.Cells.Item(80, 1).Value = 79
$ws.Cells.Item(80, 2).Value = '帝国機神ヴォルカミオン2'
$ws.Cells.Item(80, 3).Value = 2
$ws.Cells.Item(80, 3).Interior.Color = 13499135
$ws.Cells.Item(81, 1).Value = 80
$ws.Cells.Item(81, 2).Value = '先日救っていただいたドラゴンです(ノヴァコミックス)3'
$ws.Cells.Item(81, 3).Value = 3
$ws.Cells.Item(81, 3).Interior.Color = 13499135
$ws.Cells.Item(82, 1).Value = 81
$ws.Cells.Item(82, 2).Value = 'レンズの向こうの女神たち2'
$ws.Cells.Item(82, 3).Value = 2
$ws.Cells.Item(82, 3).Interior.Color = 13499135
$ws.Cells.Item(83, 1).Value = 82
$ws.Cells.Item(83, 2).Value = '水魔法なんて使えないと追放されたけど、水が万能だと気がつき水の賢者と呼ばれるまでに成長しました~今更水不足と泣きついても簡単には譲れません~3'
$ws.Cells.Item(83, 3).Value = 3
$ws.Cells.Item(83, 3).Interior.Color = 13499135
$ws.Cells.Item(84, 1).Value = 83
$ws.Cells.Item(84, 2).Value = 'ブリッツ・マジック・スケーリング@COMIC 第1話'
$ws.Cells.Item(84, 3).Value = 1
$ws.Cells.Item(84, 3).Interior.Color = 13499135
$ws.Cells.Item(85, 1).Value = 84
$ws.Cells.Item(85, 2).Value = 'え? ギルド内で唯一を極めてる俺をクビですか?@COMIC 第1話'
$ws.Cells.Item(85, 3).Value = 1
$ws.Cells.Item(85, 3).Interior.Color = 13499135
$ws.Cells.Item(86, 1).Value = 85
$ws.Cells.Item(86, 2).Value = '転生したらスライムだった件~魔物の国の歩き方~'
$ws.Cells.Item(86, 3).Value = 1
$ws.Cells.Item(86, 3).Interior.Color = 13499135
$ws.Cells.Item(87, 1).Value = 86
$ws.Cells.Item(87, 2).Value = '元騎士の辺境伯令嬢は悪魔の花嫁となる1'
$ws.Cells.Item(87, 3).Value = 1
$ws.Cells.Item(87, 3).Interior.Color = 13499135
$ws.Cells.Item(88, 1).Value = 87
$ws.Cells.Item(88, 2).Value = '恋するふたりは裏の顔がある'
$ws.Cells.Item(88, 3).Value = 2
$ws.Cells.Item(88, 3).Interior.Color = 13499135
$ws.Cells.Item(89, 1).Value = 88
$ws.Cells.Item(89, 2).Value = '直葉くんはすぐに結婚したい!'
$ws.Cells.Item(89, 3).Value = 3
$ws.Cells.Item(89, 3).Interior.Color = 13499135
$ws.Cells.Item(90, 1).Value = 89
$ws.Cells.Item(90, 2).Value = '機密少女と暗号戦争'
$ws.Cells.Item(90, 3).Value = 4
$ws.Cells.Item(91, 1).Value = 90
$ws.Cells.Item(91, 2).Value = '私と上司の内緒の事情'
$ws.Cells.Item(91, 3).Value = 11
$ws.Cells.Item(92, 1).Value = 91
$ws.Cells.Item(92, 2).Value = '灰かぶりの天使'
$ws.Cells.Item(92, 3).Value = 2
$ws.Cells.Item(92, 3).Interior.Color = 13499135
$ws.Cells.Item(93, 1).Value = 92
$ws.Cells.Item(93, 2).Value = '灰かぶりの天使'
$ws.Cells.Item(93, 3).Value = 3
$ws.Cells.Item(93, 3).Interior.Color = 13499135
$ws.Cells.Item(94, 1).Value = 93
$ws.Cells.Item(94, 2).Value = 'ドリーミングスクール'
$ws.Cells.Item(94, 3).Value = 2
$ws.Cells.Item(94, 3).Interior.Color = 13499135
$ws.Cells.Item(95, 1).Value = 94
$ws.Cells.Item(95, 2).Value = '僕の吸血姫を笑わせたい'
$ws.Cells.Item(95, 3).Value = 1
$ws.Cells.Item(95, 3).Interior.Color = 13499135
$ws.Cells.Item(96, 1).Value = 95
$ws.Cells.Item(96, 2).Value = '転生したらポンコツメイドと呼ばれていました'
$ws.Cells.Item(96, 3).Value = 2
$ws.Cells.Item(96, 3).Interior.Color = 13499135
$ws.Cells.Item(97, 1).Value = 96
$ws.Cells.Item(97, 2).Value = '3組の幼なじみ'
$ws.Cells.Item(97, 3).Value = 1
$ws.Cells.Item(97, 3).Interior.Color = 13499135
$ws.Cells.Item(98, 1).Value = 97
$ws.Cells.Item(98, 2).Value = 'お父さんが早く死にますように。2'
$ws.Cells.Item(98, 3).Value = 2
$ws.Cells.Item(98, 3).Interior.Color = 13499135
$ws.Cells.Item(99, 1).Value = 98
$ws.Cells.Item(99, 2).Value = 'ぼっち・ざ・ろっく!'
$ws.Cells.Item(99, 3).Value = 8
$ws.Cells.Item(100, 1).Value = 99
$ws.Cells.Item(100, 2).Value = 'GIANT KILLING'
$ws.Cells.Item(100, 3).Value = 68
$ws.Cells.Item(101, 1).Value = 100
$ws.Cells.Item(101, 2).Value = '外れスキル「世界図書館」による異世界の知識と始める『産業革命』 ファイアーアロー?うるせえ、こっちはライフルだ!!'
$ws.Cells.Item(101, 3).Value = 4
